$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 67.875
$ws.Range("I9").Value = 54
$ws.Range("J9").Value = 91
$ws.Range("K9").Value = 54
$ws.Range("L9").Value = 91
$ws.Range("M9").Value = 115
$ws.Range("N9").Value = -429
$ws.Range("H19").Value = 941
$ws.Range("I19").Value = 854.55554
$ws.Range("J19").Value = 1096.6
$ws.Range("K19").Value = 854.55554
$ws.Range("L19").Value = 1096.6
$ws.Range("M19").Value = -679.55554
$ws.Range("N19").Value = -1446.6
$ws.Range("H31").Value = 99
$ws.Range("I31").Value = 99
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 297
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -67
$ws.Range("H42").Value = 4999.5
$ws.Range("I42").Value = 4999
$ws.Range("J42").Value = 5000
$ws.Range("K42").Value = 14997
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = -14767
$ws.Range("N42").Value = -15460
$ws.Range("H80").Value = 2462.5
$ws.Range("I80").Value = 2416.6667
$ws.Range("J80").Value = 2600
$ws.Range("K80").Value = 7250.000100000001
$ws.Range("L80").Value = 7800
$ws.Range("M80").Value = -6252.000100000001
$ws.Range("N80").Value = -9796
$ws.Range("H83").Value = 2462.5
$ws.Range("I83").Value = 2416.6667
$ws.Range("J83").Value = 2600
$ws.Range("K83").Value = 21750.0003
$ws.Range("L83").Value = 23400
$ws.Range("M83").Value = -16758.0003
$ws.Range("N83").Value = -33384
$ws.Range("H92").Value = 688.4
$ws.Range("I92").Value = 298.25
$ws.Range("J92").Value = 2249
$ws.Range("K92").Value = 298.25
$ws.Range("L92").Value = 2249
$ws.Range("M92").Value = 949.75
$ws.Range("N92").Value = -4745
$ws.Range("H101").Value = 4317.5
$ws.Range("I101").Value = 6785
$ws.Range("J101").Value = 1850
$ws.Range("K101").Value = 20355
$ws.Range("L101").Value = 5550
$ws.Range("M101").Value = -18733
$ws.Range("N101").Value = -8794
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = ""
$ws.Range("H132").Value = 9384.5
$ws.Range("I132").Value = 10061.5
$ws.Range("J132").Value = 5999.5
$ws.Range("K132").Value = 30184.5
$ws.Range("L132").Value = 17998.5
$ws.Range("M132").Value = -27654.5
$ws.Range("N132").Value = -23058.5
$ws.Range("H138").Value = 1476.8125
$ws.Range("I138").Value = 559
$ws.Range("J138").Value = 2656.8572
$ws.Range("K138").Value = 1677
$ws.Range("L138").Value = 7970.571599999999
$ws.Range("M138").Value = 3463
$ws.Range("N138").Value = -18250.5716
$ws.Range("H141").Value = 1699.75
$ws.Range("I141").Value = 1699.75
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5099.25
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 80.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 59999.75
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 59999.75
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 59999.75
$ws.Range("N76").Value = -60675.75
$ws.Range("H79").Value = 59999.75
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 59999.75
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 59999.75
$ws.Range("N79").Value = -62339.75
$ws.Range("H95").Value = 15301.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 15301.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 15301.5
$ws.Range("N95").Value = -20793.5
$ws.Range("H122").Value = 1200.625
$ws.Range("I122").Value = 1186
$ws.Range("J122").Value = 1244.5
$ws.Range("K122").Value = 3558
$ws.Range("L122").Value = 3733.5
$ws.Range("M122").Value = -1108
$ws.Range("N122").Value = -8633.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3374.75
$ws.Range("I20").Value = 4263.3335
$ws.Range("J20").Value = 709
$ws.Range("K20").Value = 4263.3335
$ws.Range("L20").Value = 709
$ws.Range("M20").Value = -4016.3335
$ws.Range("N20").Value = -1203
$ws.Range("H22").Value = 600.8
$ws.Range("I22").Value = 667.3333
$ws.Range("J22").Value = 501
$ws.Range("K22").Value = 667.3333
$ws.Range("L22").Value = 501
$ws.Range("M22").Value = -494.3333
$ws.Range("N22").Value = -847
$ws.Range("H46").Value = 100000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 100000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 100000
$ws.Range("N46").Value = -100596
$ws.Range("H49").Value = 100000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 100000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 100000
$ws.Range("N49").Value = -100478
$ws.Range("H64").Value = 839.6667
$ws.Range("I64").Value = 609.5
$ws.Range("J64").Value = 1300
$ws.Range("K64").Value = 609.5
$ws.Range("L64").Value = 1300
$ws.Range("M64").Value = -384.5
$ws.Range("N64").Value = -1750
$ws.Range("H67").Value = 839.6667
$ws.Range("I67").Value = 609.5
$ws.Range("J67").Value = 1300
$ws.Range("K67").Value = 609.5
$ws.Range("L67").Value = 1300
$ws.Range("M67").Value = 170.5
$ws.Range("N67").Value = -2860
$ws.Range("H86").Value = 6246.727
$ws.Range("I86").Value = 4914.625
$ws.Range("J86").Value = 9799
$ws.Range("K86").Value = 4914.625
$ws.Range("L86").Value = 9799
$ws.Range("M86").Value = -3791.625
$ws.Range("N86").Value = -12045
$ws.Range("H89").Value = 6246.727
$ws.Range("I89").Value = 4914.625
$ws.Range("J89").Value = 9799
$ws.Range("K89").Value = 24573.125
$ws.Range("L89").Value = 48995
$ws.Range("M89").Value = -18957.125
$ws.Range("N89").Value = -60227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 10.333333
$ws.Range("I7").Value = 10.2
$ws.Range("J7").Value = 10.5
$ws.Range("K7").Value = 10.2
$ws.Range("L7").Value = 10.5
$ws.Range("M7").Value = 102.8
$ws.Range("N7").Value = -236.5
$ws.Range("H31").Value = 3092.375
$ws.Range("I31").Value = 1958.8
$ws.Range("J31").Value = 4981.6665
$ws.Range("K31").Value = 1958.8
$ws.Range("L31").Value = 4981.6665
$ws.Range("M31").Value = -1663.8
$ws.Range("N31").Value = -5571.6665
$ws.Range("H34").Value = 3092.375
$ws.Range("I34").Value = 1958.8
$ws.Range("J34").Value = 4981.6665
$ws.Range("K34").Value = 1958.8
$ws.Range("L34").Value = 4981.6665
$ws.Range("M34").Value = -1756.8
$ws.Range("N34").Value = -5385.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 57.636364
$ws.Range("I2").Value = 77.666664
$ws.Range("J2").Value = 33.6
$ws.Range("K2").Value = 465.999984
$ws.Range("L2").Value = 201.6
$ws.Range("M2").Value = -352.999984
$ws.Range("N2").Value = -427.6
$ws.Range("H4").Value = 173.63637
$ws.Range("I4").Value = 173.63637
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 520.9091100000001
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -408.9091100000001
$ws.Range("N4").Value = ""
$ws.Range("H11").Value = 52
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 53
$ws.Range("K11").Value = 150
$ws.Range("L11").Value = 159
$ws.Range("M11").Value = -10
$ws.Range("N11").Value = -439
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = ""
$ws.Range("H38").Value = 981.2222
$ws.Range("I38").Value = 3786.5
$ws.Range("J38").Value = 179.71428
$ws.Range("K38").Value = 11359.5
$ws.Range("L38").Value = 539.14284
$ws.Range("M38").Value = -11012.5
$ws.Range("N38").Value = -1233.14284
$ws.Range("H103").Value = 736
$ws.Range("I103").Value = 726.6667
$ws.Range("J103").Value = 750
$ws.Range("K103").Value = 2180.0001
$ws.Range("L103").Value = 2250
$ws.Range("M103").Value = -1301.0001
$ws.Range("N103").Value = -4008
$ws.Range("H107").Value = 665
$ws.Range("I107").Value = 497.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1492.5
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 427.5
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 20
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 20
$ws.Range("L2").Value = 20
$ws.Range("M2").Value = 93
$ws.Range("N2").Value = -246
$ws.Range("H70").Value = 7667.4614
$ws.Range("I70").Value = 6961
$ws.Range("J70").Value = 8109
$ws.Range("K70").Value = 6961
$ws.Range("L70").Value = 8109
$ws.Range("M70").Value = -6691
$ws.Range("N70").Value = -8649
$ws.Range("H73").Value = 7667.4614
$ws.Range("I73").Value = 6961
$ws.Range("J73").Value = 8109
$ws.Range("K73").Value = 6961
$ws.Range("L73").Value = 8109
$ws.Range("M73").Value = -6025
$ws.Range("N73").Value = -9981
$ws.Range("H80").Value = 11952.667
$ws.Range("I80").Value = 3939.3333
$ws.Range("J80").Value = 19966
$ws.Range("K80").Value = 3939.3333
$ws.Range("L80").Value = 19966
$ws.Range("M80").Value = -2941.3333
$ws.Range("N80").Value = -21962
$ws.Range("H83").Value = 11952.667
$ws.Range("I83").Value = 3939.3333
$ws.Range("J83").Value = 19966
$ws.Range("K83").Value = 19696.6665
$ws.Range("L83").Value = 99830
$ws.Range("M83").Value = -14704.6665
$ws.Range("N83").Value = -109814
$ws.Range("H102").Value = 1499.5
$ws.Range("I102").Value = 1499.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1499.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 122.5
$ws.Range("N102").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2013.8572
$ws.Range("I22").Value = 1724.25
$ws.Range("J22").Value = 2400
$ws.Range("K22").Value = 1724.25
$ws.Range("L22").Value = 2400
$ws.Range("M22").Value = -1429.25
$ws.Range("N22").Value = -2990
$ws.Range("H27").Value = 2013.8572
$ws.Range("I27").Value = 1724.25
$ws.Range("J27").Value = 2400
$ws.Range("K27").Value = 1724.25
$ws.Range("L27").Value = 2400
$ws.Range("M27").Value = -1617.25
$ws.Range("N27").Value = -2614
$ws.Range("H122").Value = 6486.1665
$ws.Range("I122").Value = 5783.4
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 17350.2
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -14900.2
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 900
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2700
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -170

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -939
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4696
$ws.Range("H122").Value = 2035.7142
$ws.Range("I122").Value = 890
$ws.Range("J122").Value = 4900
$ws.Range("K122").Value = 2670
$ws.Range("L122").Value = 14700
$ws.Range("M122").Value = -220
$ws.Range("N122").Value = -19600
$ws.Range("H123").Value = 157500
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 157500
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 157500
$ws.Range("N123").Value = -167300
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""

